$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'65.880.03"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -2.04%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'3.258.98"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -1.82%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  -0.13%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'573.37"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -0.80%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'177.62"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -4.73%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.622"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  +2.37%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.999"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -0.13%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = "'  -3.97%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'6.69"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +0.19%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = "'  -2.23%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'3.820.61"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -1.66%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("E13").Value = "'  -3.87%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'65.892.57"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -2.41%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'26.32"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -4.05%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("B16").Value = "'ShibaInu"
$ws.Range("B16").Style = "Normal"
$ws.Range("C16").Value = "'https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("C16").Style = "Normal"
$ws.Range("D16").Value = "'0.0000161"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -3.76%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("B17").Value = "'WrappedEther"
$ws.Range("B17").Style = "Normal"
$ws.Range("C17").Value = "'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("C17").Style = "Normal"
$ws.Range("D17").Value = "'3.249.86"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -2.03%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'430.77"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -3.26%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'5.52"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -3.01%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Value = "'  -3.68%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'7.35"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -4.92%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'71.67"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -3.25%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("E23").Value = "'  -0.07%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'3.409.37"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -1.38%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.502"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -2.90%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.195"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +3.51%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'0.0000111"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -6.23%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'8.83"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -2.45%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'1.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -0.08%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("E30").Value = "'  -2.63%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'22.13"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -3.51%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("E32").Value = "'  +0.13%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("E33").Value = "'  -4.13%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("E34").Value = "'  -4.19%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("E35").Value = "'  -5.61%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'157.11"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -3.65%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("E37").Value = "'  -6.78%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("B38").Value = "'Stacks"
$ws.Range("B38").Style = "Normal"
$ws.Range("C38").Value = "'https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("C38").Style = "Normal"
$ws.Range("D38").Value = "'1.78"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -3.77%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("B39").Value = "'EnergySwap"
$ws.Range("B39").Style = "Normal"
$ws.Range("C39").Value = "'https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("C39").Style = "Normal"
$ws.Range("D39").Value = "'26.26"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -3.44%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'2.759.12"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -0.73%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.771"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -2.40%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("E42").Value = "'  -4.39%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'40.19"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +0.27%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'6.02"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -3.48%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.0651"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -2.87%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'319.29"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -2.10%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = "'  -5.22%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'23.15"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -6.66%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'0.0264"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -3.29%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = "'  +1.45%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.999"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -0.04%  "
$ws.Range("E51").Style = "Normal"
